$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (F2:K2) previously held the long descriptive dataset names; switch
# them to the short internal codes already used in row 3, matching the
# corrected lookup keys.
$ws.Range("F2").Value = "rail_trans"
$ws.Range("G2").Value = "wate_trans"
$ws.Range("H2").Value = "land_trans"
$ws.Range("I2").Value = "truck_trans"
$ws.Range("J2").Value = "sea_trans"
$ws.Range("K2").Value = "powder"

# Unit conversion of powder surface density: 12 -> 0.48 (kg/m**2)
$ws.Range("K9").Value = 0.48

# Document the basis for the converted value with a threaded comment.
$excel.UserName = "Lorenzo Rinaldi"
$ws.Range("K9").AddCommentThreaded("60 micron thickness, 8000 kg/m3 densitity assumed for steel powder coating`n")
